# Update the practice-problem answers in the single table of the document.
# Each table cell is addressed explicitly (row, column) so that duplicate
# values (e.g. "39÷7=5, 4" appears twice) are replaced independently with
# their correct new values, rather than relying on a global text search.

$d = $word.ActiveDocument

$replacements = @(
    @{ Row = 1;  Col = 1; Old = "18÷9=2, 0";   New = "30÷2=15, 0" },
    @{ Row = 1;  Col = 2; Old = "39÷7=5, 4";   New = "23÷7=3, 2" },
    @{ Row = 1;  Col = 3; Old = "87÷7=12, 3";  New = "94÷8=11, 6" },
    @{ Row = 1;  Col = 4; Old = "57÷3=19, 0";  New = "19÷7=2, 5" },
    @{ Row = 1;  Col = 5; Old = "87÷6=14, 3";  New = "92÷3=30, 2" },

    @{ Row = 5;  Col = 1; Old = "73÷2=36, 1";  New = "70÷8=8, 6" },
    @{ Row = 5;  Col = 2; Old = "58÷9=6, 4";   New = "54÷3=18, 0" },
    @{ Row = 5;  Col = 3; Old = "80÷7=11, 3";  New = "55÷2=27, 1" },
    @{ Row = 5;  Col = 4; Old = "80÷8=10, 0";  New = "76÷3=25, 1" },
    @{ Row = 5;  Col = 5; Old = "39÷9=4, 3";   New = "85÷3=28, 1" },

    @{ Row = 9;  Col = 1; Old = "34÷6=5, 4";   New = "32÷4=8, 0" },
    @{ Row = 9;  Col = 2; Old = "93÷6=15, 3";  New = "23÷5=4, 3" },
    @{ Row = 9;  Col = 3; Old = "39÷7=5, 4";   New = "74÷6=12, 2" },
    @{ Row = 9;  Col = 4; Old = "44÷2=22, 0";  New = "35÷2=17, 1" },
    @{ Row = 9;  Col = 5; Old = "88÷9=9, 7";   New = "72÷6=12, 0" },

    @{ Row = 13; Col = 1; Old = "83÷6=13, 5";  New = "28÷2=14, 0" },
    @{ Row = 13; Col = 2; Old = "57÷2=28, 1";  New = "16÷3=5, 1" },
    @{ Row = 13; Col = 3; Old = "44÷9=4, 8";   New = "55÷6=9, 1" },
    @{ Row = 13; Col = 4; Old = "30÷9=3, 3";   New = "96÷7=13, 5" },
    @{ Row = 13; Col = 5; Old = "39÷4=9, 3";   New = "90÷3=30, 0" },

    @{ Row = 17; Col = 1; Old = "31÷3=10, 1";  New = "69÷6=11, 3" },
    @{ Row = 17; Col = 2; Old = "17÷2=8, 1";   New = "60÷9=6, 6" },
    @{ Row = 17; Col = 3; Old = "60÷8=7, 4";   New = "92÷6=15, 2" },
    @{ Row = 17; Col = 4; Old = "90÷2=45, 0";  New = "84÷3=28, 0" },
    @{ Row = 17; Col = 5; Old = "43÷7=6, 1";   New = "11÷6=1, 5" }
)

foreach ($item in $replacements) {
    # Re-fetch the table/cell fresh on every iteration rather than caching a
    # single handle, so a previous edit cannot leave us pointing at a stale
    # cell reference.
    $cellRange = $d.Tables.Item(1).Cell($item.Row, $item.Col).Range
    # Use wdReplaceOne (1) rather than wdReplaceAll (2): several answers in
    # this worksheet repeat the same text (e.g. "39÷7=5, 4" occurs twice),
    # and a ReplaceAll initiated from a single cell's Range would otherwise
    # replace every matching occurrence in the whole document instead of
    # just the targeted cell.
    $cellRange.Find.Execute($item.Old, $true, $false, $false, $false, $false, `
                             $true, 0, $false, $item.New, 1)
}
